# Update Name of Algo
# Apply updated KNN imputation results to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A3"   = -21.933
    "D6"   = -8.007000000000001
    "A14"  = -21.644
    "D18"  = -8.289
    "D19"  = -7.884
    "A21"  = -20.866
    "B22"  = 6.881
    "A23"  = -20.502
    "B24"  = 5.342000000000001
    "A25"  = -21.589
    "A26"  = -21.367
    "B28"  = 5.468000000000001
    "A29"  = -21.577
    "B36"  = 6.19
    "D44"  = -7.869000000000002
    "B45"  = 5.646999999999999
    "D47"  = -7.662000000000001
    "B48"  = 6.238
    "B49"  = 6.237
    "D51"  = -8.231
    "B52"  = 5.58
    "A53"  = -20.72
    "B53"  = 8.151
    "B54"  = 5.144000000000001
    "D55"  = -8.288999999999998
    "A57"  = -22.16
    "D57"  = -8.074999999999999
    "A59"  = -22.282
    "D64"  = -7.407000000000001
    "A69"  = -21.441
    "B70"  = 5.511000000000001
    "A79"  = -21.348
    "D80"  = -8.166
    "A83"  = -21.963
    "B86"  = 4.936999999999999
    "B87"  = 4.537000000000001
    "B89"  = 4.983000000000001
    "A91"  = -21.062
    "D92"  = -7.355000000000001
    "A93"  = -21.533
    "D94"  = -7.703999999999999
    "D96"  = -7.464
    "B101" = 5.731999999999999
    "D101" = -7.540000000000001
    "A103" = -22.055
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
